# Add 2022-Q1 sheet (feat: add 2022-Q1 data)
#
# Workbook currently has two sheets: "2021-Q4" and "总计" (totals).
# We insert a new "2022-Q1" sheet (same layout/shape as "2021-Q4") right
# after "2021-Q4", fill it with the new quarter's fund-holding data, and
# update the "总计" sheet with a new leading row summarising 2022-Q1
# (pushing the existing 2021-Q4 summary row down one row).

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) Build the "2022-Q1" sheet by duplicating "2021-Q4" (so it starts out
#    with the exact same header row / column layout / cell styles), then
#    placing it immediately after "2021-Q4".
# ---------------------------------------------------------------------
$q4Sheet.Copy([System.Reflection.Missing]::Value, $q4Sheet)
$q1Sheet = $wb.Worksheets.Item(2)
$q1Sheet.Name = "2022-Q1"

# The duplicated sheet carries 11 data rows (rows 2-12); 2022-Q1 only has
# 4 holdings, so drop the surplus rows 6-12.
$q1Sheet.Range("A6:H12").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2) Write the 2022-Q1 holdings. Columns B,C,D,E,F,G hold numeric-looking
#    text (fund codes with leading zeros, percentages, etc.) so mark them
#    as text before assigning, same as the source data's formatting.
# ---------------------------------------------------------------------
function Set-FundRow {
    param($sheet, $row, $idx, $code, $name, $size, $position, $pct, $value, $rank)

    $sheet.Range("B$row`:G$row").NumberFormat = "@"

    $sheet.Range("A$row").Value = $idx
    $sheet.Range("B$row").Value = $code
    $sheet.Range("C$row").Value = $name
    $sheet.Range("D$row").Value = $size
    $sheet.Range("E$row").Value = $position
    $sheet.Range("F$row").Value = $pct
    $sheet.Range("G$row").Value = $value
    $sheet.Range("H$row").Value = $rank
}

Set-FundRow $q1Sheet 2 0 "000800" "华商未来主题混合"             "4.31" "84.71" "3.62" "0.1560" 7
Set-FundRow $q1Sheet 3 1 "010487" "中银顺盈回报一年持有期混合"     "1.23" "38.83" "0.65" "0.0080" 10
Set-FundRow $q1Sheet 4 2 "001914" "中信建投聚利混合A"             "0.13" "39.07" "2.10" "0.0027" 4
Set-FundRow $q1Sheet 5 3 "000041" "华夏全球精选股票(QDII)"        "0.02" "39.07" "2.10" "0.0004" 4

$q1Sheet.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: insert a new 2022-Q1 summary row
#    above the existing 2021-Q4 row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push the current row 2 (2021-Q4 totals) down to row 3.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 11
$totalSheet.Range("D3").Value = 4.94

# Match the bold/centered/bordered look used on the other index cells
# (column A) in this workbook.
$totalSheet.Range("A3").Font.Bold = $true
$totalSheet.Range("A3").Borders.LineStyle = 1
$totalSheet.Range("A3").HorizontalAlignment = -4108
$totalSheet.Range("A3").VerticalAlignment = -4160

# Write the new 2022-Q1 summary into row 2.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.17

$totalSheet.Range("A1").Select() | Out-Null
